# Auto-generated edit script applying the crypto price/volume update
# described by the commit "Updated cryptos list on Sat Dec  2 20:36:40 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainValue($addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-TextValue($addr, $val) {
    # Force the cell to keep a literal text representation (e.g. "1.00", "16.19")
    # instead of Excel silently auto-converting the numeric-looking string to a
    # floating point number. Reset the style back to Normal afterwards so no
    # stray formatting is left behind on the cell.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-PlainValue "D2" '39.377.45'
Set-PlainValue "E2" '  +1.38%  '

# Row 3
Set-PlainValue "D3" '2.157.36'
Set-PlainValue "E3" '  +3.11%  '

# Row 4
Set-TextValue "D4" '1.00'
Set-PlainValue "E4" '  -0.03%  '

# Row 5
Set-TextValue "D5" '229.33'
Set-PlainValue "E5" '  +0.23%  '

# Row 6
Set-TextValue "D6" '0.622'
Set-PlainValue "E6" '  +1.26%  '

# Row 7
Set-TextValue "D7" '63.17'
Set-PlainValue "E7" '  +4.11%  '

# Row 8
Set-PlainValue "E8" '  +0.06%  '

# Row 9
Set-PlainValue "E9" '  +2.26%  '

# Row 10
Set-TextValue "D10" '0.0859'
Set-PlainValue "E10" '  +2.56%  '

# Row 11
Set-PlainValue "E11" '  +0.32%  '

# Row 12
Set-TextValue "D12" '16.19'
Set-PlainValue "E12" '  +8.06%  '

# Row 13
Set-PlainValue "D13" '2.479.87'
Set-PlainValue "E13" '  +3.26%  '

# Row 14
Set-TextValue "D14" '22.28'
Set-PlainValue "E14" '  +2.01%  '

# Row 15
Set-TextValue "D15" '0.817'
Set-PlainValue "E15" '  +2.57%  '

# Row 16
Set-PlainValue "E16" '  +1.49%  '

# Row 17
Set-PlainValue "D17" '2.155.97'
Set-PlainValue "E17" '  +2.92%  '

# Row 18
Set-PlainValue "D18" '39.564.60'
Set-PlainValue "E18" '  +2.17%  '

# Row 19
Set-TextValue "D19" '72.44'
Set-PlainValue "E19" '  +0.84%  '

# Row 20
Set-PlainValue "E20" '  +1.87%  '

# Row 21
Set-PlainValue "D21" '0.0₃0853'
Set-PlainValue "E21" '  +1.75%  '

# Row 22
Set-TextValue "D22" '229.09'
Set-PlainValue "E22" '  +0.85%  '

# Row 24
Set-TextValue "D24" '2.42'
Set-PlainValue "E24" '  +1.56%  '

# Row 25
Set-TextValue "D25" '2.37'
Set-PlainValue "E25" '  +1.46%  '

# Row 26
Set-TextValue "D26" '9.69'
Set-PlainValue "E26" '  +2.28%  '

# Row 27
Set-TextValue "D27" '173.10'
Set-PlainValue "E27" '  +1.15%  '

# Row 28
Set-PlainValue "E28" '  -1.01%  '

# Row 29
Set-PlainValue "B29" 'ImmutableX'
Set-PlainValue "C29" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D29" '1.42'
Set-PlainValue "E29" '  -2.25%  '

# Row 30
Set-PlainValue "B30" 'EthereumClassic'
Set-PlainValue "C30" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D30" '19.65'
Set-PlainValue "E30" '  +2.36%  '

# Row 31
Set-PlainValue "E31" '  +7.79%  '

# Row 32
Set-PlainValue "E32" '  +1.31%  '

# Row 33
Set-TextValue "D33" '4.66'
Set-PlainValue "E33" '  +3.44%  '

# Row 34
Set-PlainValue "E34" '  +3.00%  '

# Row 35
Set-PlainValue "E35" '  +11.82%  '

# Row 36
Set-TextValue "D36" '0.0622'
Set-PlainValue "E36" '  +1.51%  '

# Row 37
Set-PlainValue "E37" '  +1.23%  '

# Row 38
Set-TextValue "D38" '3.58'
Set-PlainValue "E38" '  +0.13%  '

# Row 39
Set-PlainValue "E39" '  +0.08%  '

# Row 40
Set-PlainValue "E40" '  +0.48%  '

# Row 41
Set-PlainValue "E41" '  +3.02%  '

# Row 42
Set-TextValue "D42" '103.33'
Set-PlainValue "E42" '  +2.47%  '

# Row 43
Set-PlainValue "D43" '1.540.72'
Set-PlainValue "E43" '  -0.09%  '

# Row 44
Set-PlainValue "E44" '  +6.54%  '

# Row 45
Set-PlainValue "E45" '  +7.01%  '

# Row 46
Set-PlainValue "E46" '  +0.52%  '

# Row 47
Set-TextValue "D47" '2.80'
Set-PlainValue "E47" '  -0.52%  '

# Row 48
Set-TextValue "D48" '7.81'
Set-PlainValue "E48" '  +2.02%  '

# Row 49
Set-PlainValue "E49" '  +1.56%  '

# Row 50
Set-PlainValue "D50" '2.363.61'
Set-PlainValue "E50" '  +3.30%  '

# Row 51
Set-PlainValue "E51" '  -0.11%  '

